$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "{Diagram}") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Delete the whole paragraph, including its trailing paragraph mark,
    # so the empty paragraph before it merges with the heading that follows.
    $target.Range.Delete()
}
